$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 30.06040833333333
$ws.Range("H2").Value = 90.181225
$ws.Range("I2").Value = 0.1875845602414817
$ws.Range("J2").Value = 0.1875845602414817
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 52.80829433333333
$ws.Range("N2").Value = 158.424883
$ws.Range("O2").Value = 0.1724060238174878
$ws.Range("P2").Value = 0.1724060238174878
$ws.Range("Q2").Value = 1587.438891046853
$ws.Range("R2").Value = 14286.95001942167
$ws.Range("S2").Value = 0.03234070816078587
$ws.Range("T2").Value = 0.03234070816078587

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 30.06040833333333
$ws.Range("H3").Value = 90.181225
$ws.Range("I3").Value = 0.1875845602414817
$ws.Range("J3").Value = 0.1875845602414817
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 81.07766966666667
$ws.Range("N3").Value = 243.233009
$ws.Range("O3").Value = 0.2646985445010758
$ws.Range("P3").Value = 0.2646985445010758
$ws.Range("Q3").Value = 2437.227856895114
$ws.Range("R3").Value = 21935.05071205602
$ws.Range("S3").Value = 0.04965336006679457
$ws.Range("T3").Value = 0.04965336006679457

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 30.06040833333333
$ws.Range("H4").Value = 90.181225
$ws.Range("I4").Value = 0.1875845602414817
$ws.Range("J4").Value = 0.1875845602414817
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 172.4159456666667
$ws.Range("N4").Value = 517.247837
$ws.Range("O4").Value = 0.5628954316814363
$ws.Range("P4").Value = 0.5628954316814364
$ws.Range("Q4").Value = 5182.893729917813
$ws.Range("R4").Value = 46646.04356926033
$ws.Range("S4").Value = 0.1055904920139012
$ws.Range("T4").Value = 0.1055904920139012

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 105.79319
$ws.Range("H5").Value = 317.3795699999999
$ws.Range("I5").Value = 0.6601762957653385
$ws.Range("J5").Value = 0.6601762957653385
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 52.80829433333333
$ws.Range("N5").Value = 158.424883
$ws.Range("O5").Value = 0.1724060238174878
$ws.Range("P5").Value = 0.1724060238174878
$ws.Range("Q5").Value = 5586.757915982255
$ws.Range("R5").Value = 50280.8212438403
$ws.Range("S5").Value = 0.1138183701714598
$ws.Range("T5").Value = 0.1138183701714598

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 105.79319
$ws.Range("H6").Value = 317.3795699999999
$ws.Range("I6").Value = 0.6601762957653385
$ws.Range("J6").Value = 0.6601762957653385
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 81.07766966666667
$ws.Range("N6").Value = 243.233009
$ws.Range("O6").Value = 0.2646985445010758
$ws.Range("P6").Value = 0.2646985445010758
$ws.Range("Q6").Value = 8577.465311802902
$ws.Range("R6").Value = 77197.18780622611
$ws.Range("S6").Value = 0.1747477046031968
$ws.Range("T6").Value = 0.1747477046031968

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 105.79319
$ws.Range("H7").Value = 317.3795699999999
$ws.Range("I7").Value = 0.6601762957653385
$ws.Range("J7").Value = 0.6601762957653385
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 172.4159456666667
$ws.Range("N7").Value = 517.247837
$ws.Range("O7").Value = 0.5628954316814363
$ws.Range("P7").Value = 0.5628954316814364
$ws.Range("Q7").Value = 18240.43289894334
$ws.Range("R7").Value = 164163.8960904901
$ws.Range("S7").Value = 0.3716102209906818
$ws.Range("T7").Value = 0.3716102209906819

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 24.39630866666667
$ws.Range("H8").Value = 73.18892600000001
$ws.Range("I8").Value = 0.1522391439931798
$ws.Range("J8").Value = 0.1522391439931798
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 52.80829433333333
$ws.Range("N8").Value = 158.424883
$ws.Range("O8").Value = 0.1724060238174878
$ws.Range("P8").Value = 0.1724060238174878
$ws.Range("Q8").Value = 1288.327448716184
$ws.Range("R8").Value = 11594.94703844566
$ws.Range("S8").Value = 0.02624694548524211
$ws.Range("T8").Value = 0.02624694548524211

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 24.39630866666667
$ws.Range("H9").Value = 73.18892600000001
$ws.Range("I9").Value = 0.1522391439931798
$ws.Range("J9").Value = 0.1522391439931798
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 81.07766966666667
$ws.Range("N9").Value = 243.233009
$ws.Range("O9").Value = 0.2646985445010758
$ws.Range("P9").Value = 0.2646985445010758
$ws.Range("Q9").Value = 1977.995855162037
$ws.Range("R9").Value = 17801.96269645833
$ws.Range("S9").Value = 0.04029747983108439
$ws.Range("T9").Value = 0.04029747983108439

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 24.39630866666667
$ws.Range("H10").Value = 73.18892600000001
$ws.Range("I10").Value = 0.1522391439931798
$ws.Range("J10").Value = 0.1522391439931798
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 172.4159456666667
$ws.Range("N10").Value = 517.247837
$ws.Range("O10").Value = 0.5628954316814363
$ws.Range("P10").Value = 0.5628954316814364
$ws.Range("Q10").Value = 4206.312629539229
$ws.Range("R10").Value = 37856.81366585306
$ws.Range("S10").Value = 0.08569471867685328
$ws.Range("T10").Value = 0.0856947186768533

